$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The scraper script re-ran: the page-number ("NA") that used to sit on
# the 2025-05-26 row no longer applies there -- that row's "Numéro de
# page" cell (C118) is cleared out...
$ws.Cells.Item(118, 3).Value = ""

# ...and a brand new row is appended for 2025-05-27, carrying the same
# "Terme"/"NA"/1 values row 118 used to have.
$ws.Cells.Item(119, 1).NumberFormat = "@"
$ws.Cells.Item(119, 1).Value = "2025-05-27"
# Match the plain (unstyled) look of the other data rows instead of the
# text-format style NumberFormat just attached to this cell.
$ws.Cells.Item(119, 1).Style = $ws.Cells.Item(118, 1).Style

$ws.Cells.Item(119, 2).Value = "Rien ne nous concerne aujourd'hui !"
$ws.Cells.Item(119, 3).Value = "NA"
$ws.Cells.Item(119, 4).Value = 1
